$d = $word.ActiveDocument

# Remove every paragraph but the first one (keeps the sectPr/body intact,
# merges away all the blank/extra paragraphs that followed "TEST").
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $d.Paragraphs($i).Range.Delete()
}

# Replace the remaining paragraph's text "TEST" with "Deep text".
$d.Paragraphs(1).Range.Text = "Deep text"

# Normal style: suppress automatic hyphenation (adds <w:suppressAutoHyphens/>
# right after <w:widowControl/> in the style's paragraph properties).
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.Hyphenation = $false
